# Update sending/target cluster labels and recomputed TPM-derived statistics
# for the Ccl4-Ccr5 NATMI edge table (rows 2-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03814
$ws.Range("H2").Value = 0.11442
$ws.Range("I2").Value = 0.0004360684493923871
$ws.Range("J2").Value = 0.0004360684493923871
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.007957
$ws.Range("N2").Value = 0.023871
$ws.Range("O2").Value = 0.0002448939493579708
$ws.Range("P2").Value = 0.0002448939493579708
$ws.Range("Q2").Value = 0.00030347998
$ws.Range("R2").Value = 0.00273131982
$ws.Range("S2").Value = 0.0000001067905247621081
$ws.Range("T2").Value = 0.0000001067905247621081

$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03814
$ws.Range("H3").Value = 0.11442
$ws.Range("I3").Value = 0.0004360684493923871
$ws.Range("J3").Value = 0.0004360684493923871
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1055696666666667
$ws.Range("N3").Value = 0.316709
$ws.Range("O3").Value = 0.003249135679578298
$ws.Range("P3").Value = 0.003249135679578299
$ws.Range("Q3").Value = 0.004026427086666667
$ws.Range("R3").Value = 0.03623784378
$ws.Range("S3").Value = 0.000001416845557659189
$ws.Range("T3").Value = 0.000001416845557659189

$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03814
$ws.Range("H4").Value = 0.11442
$ws.Range("I4").Value = 0.0004360684493923871
$ws.Range("J4").Value = 0.0004360684493923871
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02449766666666667
$ws.Range("N4").Value = 0.073493
$ws.Range("O4").Value = 0.0007539688752111494
$ws.Range("P4").Value = 0.0007539688752111494
$ws.Range("Q4").Value = 0.0009343410066666667
$ws.Range("R4").Value = 0.00840906906
$ws.Range("S4").Value = 0.0000003287820383034481
$ws.Range("T4").Value = 0.0000003287820383034481

$ws.Range("A5").Value = "ECs"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03814
$ws.Range("H5").Value = 0.11442
$ws.Range("I5").Value = 0.0004360684493923871
$ws.Range("J5").Value = 0.0004360684493923871
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 32.353591
$ws.Range("N5").Value = 97.060773
$ws.Range("O5").Value = 0.9957520014958525
$ws.Range("P5").Value = 0.9957520014958525
$ws.Range("Q5").Value = 1.23396596074
$ws.Range("R5").Value = 11.10569364666
$ws.Range("S5").Value = 0.0004342160312716624
$ws.Range("T5").Value = 0.0004342160312716623

$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2972863333333333
$ws.Range("H6").Value = 0.891859
$ws.Range("I6").Value = 0.003398982443686812
$ws.Range("J6").Value = 0.003398982443686811
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.007957
$ws.Range("N6").Value = 0.023871
$ws.Range("O6").Value = 0.0002448939493579708
$ws.Range("P6").Value = 0.0002448939493579708
$ws.Range("Q6").Value = 0.002365507354333333
$ws.Range("R6").Value = 0.021289566189
$ws.Range("S6").Value = 0.0000008323902344328697
$ws.Range("T6").Value = 0.0000008323902344328696

$ws.Range("A7").Value = "FAPs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.2972863333333333
$ws.Range("H7").Value = 0.891859
$ws.Range("I7").Value = 0.003398982443686812
$ws.Range("J7").Value = 0.003398982443686811
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1055696666666667
$ws.Range("N7").Value = 0.316709
$ws.Range("O7").Value = 0.003249135679578298
$ws.Range("P7").Value = 0.003249135679578299
$ws.Range("Q7").Value = 0.03138441911455556
$ws.Range("R7").Value = 0.282459772031
$ws.Range("S7").Value = 0.00001104375513204305
$ws.Range("T7").Value = 0.00001104375513204305

$ws.Range("A8").Value = "FAPs"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2972863333333333
$ws.Range("H8").Value = 0.891859
$ws.Range("I8").Value = 0.003398982443686812
$ws.Range("J8").Value = 0.003398982443686811
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02449766666666667
$ws.Range("N8").Value = 0.073493
$ws.Range("O8").Value = 0.0007539688752111494
$ws.Range("P8").Value = 0.0007539688752111494
$ws.Range("Q8").Value = 0.007282821498555555
$ws.Range("R8").Value = 0.065545393487
$ws.Range("S8").Value = 0.000002562726969928989
$ws.Range("T8").Value = 0.000002562726969928989

$ws.Range("A9").Value = "FAPs"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2972863333333333
$ws.Range("H9").Value = 0.891859
$ws.Range("I9").Value = 0.003398982443686812
$ws.Range("J9").Value = 0.003398982443686811
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 32.353591
$ws.Range("N9").Value = 97.060773
$ws.Range("O9").Value = 0.9957520014958525
$ws.Range("P9").Value = 0.9957520014958525
$ws.Range("Q9").Value = 9.618280438556333
$ws.Range("R9").Value = 86.56452394700699
$ws.Range("S9").Value = 0.003384543571350406
$ws.Range("T9").Value = 0.003384543571350406

$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 87.12790666666666
$ws.Range("H10").Value = 261.38372
$ws.Range("I10").Value = 0.9961649491069209
$ws.Range("J10").Value = 0.9961649491069208
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.007957
$ws.Range("N10").Value = 0.023871
$ws.Range("O10").Value = 0.0002448939493579708
$ws.Range("P10").Value = 0.0002448939493579708
$ws.Range("Q10").Value = 0.6932767533466667
$ws.Range("R10").Value = 6.23949078012
$ws.Range("S10").Value = 0.0002439547685987758
$ws.Range("T10").Value = 0.0002439547685987758

$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 87.12790666666666
$ws.Range("H11").Value = 261.38372
$ws.Range("I11").Value = 0.9961649491069209
$ws.Range("J11").Value = 0.9961649491069208
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.1055696666666667
$ws.Range("N11").Value = 0.316709
$ws.Range("O11").Value = 0.003249135679578298
$ws.Range("P11").Value = 0.003249135679578299
$ws.Range("Q11").Value = 9.198064064164445
$ws.Range("R11").Value = 82.78257657748
$ws.Range("S11").Value = 0.003236675078888596
$ws.Range("T11").Value = 0.003236675078888597

$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 87.12790666666666
$ws.Range("H12").Value = 261.38372
$ws.Range("I12").Value = 0.9961649491069209
$ws.Range("J12").Value = 0.9961649491069208
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02449766666666667
$ws.Range("N12").Value = 0.073493
$ws.Range("O12").Value = 0.0007539688752111494
$ws.Range("P12").Value = 0.0007539688752111494
$ws.Range("Q12").Value = 2.134430414884444
$ws.Range("R12").Value = 19.20987373396
$ws.Range("S12").Value = 0.000751077366202917
$ws.Range("T12").Value = 0.0007510773662029169

$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 87.12790666666666
$ws.Range("H13").Value = 261.38372
$ws.Range("I13").Value = 0.9961649491069209
$ws.Range("J13").Value = 0.9961649491069208
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 32.353591
$ws.Range("N13").Value = 97.060773
$ws.Range("O13").Value = 0.9957520014958525
$ws.Range("P13").Value = 0.9957520014958525
$ws.Range("Q13").Value = 2818.900656979506
$ws.Range("R13").Value = 25370.10591281556
$ws.Range("S13").Value = 0.9919332418932305
$ws.Range("T13").Value = 0.9919332418932304

